# Adds a new "TransactionType16" column to the RawEarnings sheet, inserted
# immediately before the old column AB (ApprenticeshipContractType), with a
# value of 0 for every existing data row. Mirrors the author's commit
# "Added transaction type 16 functionality".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RawEarnings")

# Insert a new blank column at AB (28th column), pushing the existing
# ApprenticeshipContractType / LearnRefNumber columns one place to the right.
$ws.Columns.Item(28).Insert() | Out-Null

# Header for the newly inserted column.
$ws.Range("AB1").Value = "TransactionType16"

# Populate the new column with 0 for every data row (rows 2-12).
$ws.Range("AB2:AB12").Value = 0

# Make RawEarnings the active sheet/tab with the new header cell selected,
# scrolled so column V is the left-most visible column - matching the
# author's saved view state.
$ws.Select() | Out-Null
$ws.Range("AB1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 22 | Out-Null
